# Update with Correct Forecast output
$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B ("ASIN") and shift the rest right.
$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Fix up week labels in column A (W01 -> W1, etc.) and populate the new
# Week_Start_Date column (B) with the week's start date, stored as text.
$ws.Range("A2:B17").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "W1"
$ws.Cells.Item(2, 2).Value = "2025-01-05"
$ws.Cells.Item(3, 1).Value = "W2"
$ws.Cells.Item(3, 2).Value = "2025-01-12"
$ws.Cells.Item(4, 1).Value = "W3"
$ws.Cells.Item(4, 2).Value = "2025-01-19"
$ws.Cells.Item(5, 1).Value = "W4"
$ws.Cells.Item(5, 2).Value = "2025-01-26"
$ws.Cells.Item(6, 1).Value = "W5"
$ws.Cells.Item(6, 2).Value = "2025-02-02"
$ws.Cells.Item(7, 1).Value = "W6"
$ws.Cells.Item(7, 2).Value = "2025-02-09"
$ws.Cells.Item(8, 1).Value = "W7"
$ws.Cells.Item(8, 2).Value = "2025-02-16"
$ws.Cells.Item(9, 1).Value = "W8"
$ws.Cells.Item(9, 2).Value = "2025-02-23"
$ws.Cells.Item(10, 1).Value = "W9"
$ws.Cells.Item(10, 2).Value = "2025-03-02"
$ws.Cells.Item(11, 1).Value = "W10"
$ws.Cells.Item(11, 2).Value = "2025-03-09"
$ws.Cells.Item(12, 1).Value = "W11"
$ws.Cells.Item(12, 2).Value = "2025-03-16"
$ws.Cells.Item(13, 1).Value = "W12"
$ws.Cells.Item(13, 2).Value = "2025-03-23"
$ws.Cells.Item(14, 1).Value = "W13"
$ws.Cells.Item(14, 2).Value = "2025-03-30"
$ws.Cells.Item(15, 1).Value = "W14"
$ws.Cells.Item(15, 2).Value = "2025-04-06"
$ws.Cells.Item(16, 1).Value = "W15"
$ws.Cells.Item(16, 2).Value = "2025-04-13"
$ws.Cells.Item(17, 1).Value = "W16"
$ws.Cells.Item(17, 2).Value = "2025-04-20"

# The "is_holiday_week" column (now J, after the insert) should hold real
# boolean values instead of 0/1 numbers.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}

# --- Sheet "Summary" ---
$sum = $wb.Worksheets.Item("Summary")
$sum.Range("B9").NumberFormat = "@"
$sum.Range("B9").Value = "1024"
$sum.Range("B10").NumberFormat = "@"
$sum.Range("B10").Value = "491"
